$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.115.94"
$ws.Range("E2").Value = "  -0.71%  "

$ws.Range("D3").Value = "1.856.32"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'233.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.30%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").Value = "'0.4662"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.95%  "

$ws.Range("D8").Value = "'0.2813"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.36%  "

$ws.Range("D9").Value = "'0.06454"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.29%  "

$ws.Range("D10").Value = "'21.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.31%  "

$ws.Range("D11").Value = "'0.07691"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.14%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.882.86"
$ws.Range("E12").Value = "  +0.83%  "

$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "'93.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.64%  "

$ws.Range("D14").Value = "'0.6811"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "

$ws.Range("D15").Value = "'5.039"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.93%  "

$ws.Range("D16").Value = "'263.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.49%  "

$ws.Range("D17").Value = "30.107.08"
$ws.Range("E17").Value = "  -0.70%  "

$ws.Range("D18").Value = "'13.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.83%  "

$ws.Range("D19").Value = "'0.000007620"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.62%  "

$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").Value = "2.134.91"
$ws.Range("E21").Value = "  +1.09%  "

$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").Value = "'5.136"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.70%  "

$ws.Range("D24").Value = "'6.076"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.23%  "

$ws.Range("D25").Value = "'9.293"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("D26").Value = "'165.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.27%  "

$ws.Range("D27").Value = "'18.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.36%  "

$ws.Range("D28").Value = "'1.895"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.39%  "

$ws.Range("D29").Value = "'1.364"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.33%  "

$ws.Range("D30").Value = "'0.09804"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.27%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.452"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.11%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.237"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.87%  "

$ws.Range("D33").Value = "'3.972"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.45%  "

$ws.Range("D34").Value = "'0.04655"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.27%  "

$ws.Range("D35").Value = "'1.118"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.27%  "

$ws.Range("D36").Value = "'0.6904"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.64%  "

$ws.Range("D37").Value = "'2.715"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").Value = "'0.01830"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.76%  "

$ws.Range("D39").Value = "'2.739"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.04%  "

$ws.Range("D40").Value = "'6.293"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.93%  "

$ws.Range("D41").Value = "'71.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.27%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'0.9996"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'1.895"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.31%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'102.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.18%  "

$ws.Range("D45").Value = "'0.8272"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.61%  "

$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Value = "'0.4062"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.71%  "

$ws.Range("D47").Value = "'947.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.50%  "

$ws.Range("D48").Value = "'6.938"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.78%  "

$ws.Range("D49").Value = "'8.936"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.99%  "

$ws.Range("D50").Value = "'34.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.67%  "

$ws.Range("D51").Value = "'0.05588"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.59%  "
